$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44421
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15400
$ws.Range("P2").Value = 513

# Row 3
$ws.Range("D3").Value = 44376
$ws.Range("K3").Value = 18000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 18000
$ws.Range("O3").Value = 'Provincia de Limarí'
$ws.Range("P3").Value = 600

# Row 4
$ws.Range("D4").Value = 44425
$ws.Range("J4").Value = 35
$ws.Range("O4").Value = 'Provincia de Limarí'

# Row 5
$ws.Range("D5").Value = 44453
$ws.Range("J5").Value = 50

# Row 6
$ws.Range("D6").Value = 44467
$ws.Range("J6").Value = 35
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 12000
$ws.Range("M6").Value = 12000
$ws.Range("P6").Value = 400

# Row 7
$ws.Range("D7").Value = 44435
$ws.Range("J7").Value = 25
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 14000
$ws.Range("P7").Value = 467

# Row 8
$ws.Range("D8").Value = 44435
$ws.Range("O8").Value = 'Provincia del Elquí'

# Row 9
$ws.Range("D9").Value = 44449
$ws.Range("K9").Value = 12000
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 12000
$ws.Range("P9").Value = 400

# Row 10
$ws.Range("D10").Value = 44460
$ws.Range("J10").Value = 45
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 13000
$ws.Range("M10").Value = 13000
$ws.Range("P10").Value = 433

# Row 11
$ws.Range("D11").Value = 44841
$ws.Range("J11").Value = 45

# Row 12
$ws.Range("D12").Value = 44841
$ws.Range("I12").Value = 'Segunda'
$ws.Range("J12").Value = 45
$ws.Range("K12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 10000
$ws.Range("N12").Value = '$/caja 40 unidades'
$ws.Range("P12").Value = 250
$ws.Range("Q12").Value = 40

# Row 13
$ws.Range("D13").Value = 44474
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("P13").Value = 333

# Row 14
$ws.Range("D14").Value = 44418
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 30
$ws.Range("K14").Value = 15000
$ws.Range("L14").Value = 15000
$ws.Range("M14").Value = 15000
$ws.Range("N14").Value = '$/caja 30 unidades'
$ws.Range("P14").Value = 500
$ws.Range("Q14").Value = 30

# Row 15
$ws.Range("D15").Value = 44446
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 14000
$ws.Range("M15").Value = 14000
$ws.Range("P15").Value = 467

# Row 16
$ws.Range("D16").Value = 44432
$ws.Range("J16").Value = 25
$ws.Range("O16").Value = 'Provincia del Elquí'
